# DC COVID-19 Data workbook update — "14 May data update"
# Appends one new day-of-data column (serial 43964) to every sheet and
# applies a couple of historical-data corrections that shipped in the
# same commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overal Stats" — new column BR (70), header/date row = 1,
# percentage row = 15. Previous (reference) column is BQ (69).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overal Stats")
$ws1NewCol = 70
$ws1PrevCol = 69

$sheet1Data = @{
  1 = 43964
  3 = 32999
  4 = 6736
  5 = 358
  6 = 966
  8 = 345
  9 = 74
  10 = 440
  11 = 264
  12 = 176
  13 = 397
  14 = 1895
  15 = 0.76
  19 = 104
  20 = 32
  21 = 72
  22 = 108
  23 = 140
  24 = 1075
  27 = 113
  28 = 35
  29 = 78
  30 = 56
  31 = 91
  32 = 1036
  35 = 76
  36 = 26
  37 = 49
  38 = 12
  39 = 38
  40 = 245
  41 = 1
  44 = 182
  45 = 38
  46 = 138
  47 = 755
  48 = 793
  49 = 534
  50 = 1
  53 = 27
  54 = 14
  55 = 12
  56 = 17
  57 = 31
  58 = 122
  59 = 1
  61 = 10
  62 = 4
  63 = 6
  64 = 45
  65 = 49
  66 = 0
  67 = 0
  70 = 269
  71 = 330
  72 = 312
  73 = 15
  75 = 84
  76 = 19
  77 = 103
  78 = 157
  79 = 1
  81 = 79
  82 = 56
  83 = 142
  84 = 5
  85 = 13
  87 = 158
  88 = 27
  89 = 22
}

foreach ($r in $sheet1Data.Keys) {
  $cell = $ws1.Cells.Item($r, $ws1NewCol)
  $cell.Value = $sheet1Data[$r]
  $cell.NumberFormat = $ws1.Cells.Item($r, $ws1PrevCol).NumberFormat
}

# Historical correction: row 3 (People Tested Overall), column BP (68) —
# 30050 -> 31050
$ws1.Cells.Item(3, 68).Value = 31050

# ---------------------------------------------------------------------
# Sheet 2: "Total Cases by Ward" — new column AS (45), header row = 2.
# Previous column is AR (44).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Total Cases by Ward")
$ws2NewCol = 45
$ws2PrevCol = 44

$sheet2Data = @{
  2 = 43964
  3 = 908
  4 = 400
  5 = 362
  6 = 1377
  7 = 1054
  8 = 624
  9 = 979
  10 = 937
  11 = 95
}

foreach ($r in $sheet2Data.Keys) {
  $cell = $ws2.Cells.Item($r, $ws2NewCol)
  $cell.Value = $sheet2Data[$r]
  $cell.NumberFormat = $ws2.Cells.Item($r, $ws2PrevCol).NumberFormat
}

# Historical correction: column AF (32), rows 3-11 — revised 4/30 figures
$sheet2Corrections = @{
  3 = 593
  4 = 317
  5 = 294
  6 = 849
  7 = 699
  8 = 468
  9 = 701
  10 = 677
  11 = 60
}
foreach ($r in $sheet2Corrections.Keys) {
  $ws2.Cells.Item($r, 32).Value = $sheet2Corrections[$r]
}

# ---------------------------------------------------------------------
# Sheet 3: "Total Cases by Race" — new column AN (40), header row = 2.
# Previous column is AM (39).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Total Cases by Race")
$ws3NewCol = 40
$ws3PrevCol = 39

$sheet3Data = @{
  2 = 43964
  4 = 6736
  5 = 795
  6 = 1076
  7 = 3179
  8 = 93
  9 = 18
  10 = 16
  11 = 1499
  12 = 42
  14 = 1259
  15 = 1573
  16 = 3881
  17 = 23
}

foreach ($r in $sheet3Data.Keys) {
  $cell = $ws3.Cells.Item($r, $ws3NewCol)
  $cell.Value = $sheet3Data[$r]
  $cell.NumberFormat = $ws3.Cells.Item($r, $ws3PrevCol).NumberFormat
}

# ---------------------------------------------------------------------
# Sheet 4: "Lives Lost by Race" — new column AN (40), header row = 1.
# Previous column is AM (39).
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Lives Lost by Race")
$ws4NewCol = 40
$ws4PrevCol = 39

$sheet4Data = @{
  1 = 43964
  3 = 358
  4 = 5
  5 = 277
  6 = 36
  7 = 38
  8 = 2
}

foreach ($r in $sheet4Data.Keys) {
  $cell = $ws4.Cells.Item($r, $ws4NewCol)
  $cell.Value = $sheet4Data[$r]
  $cell.NumberFormat = $ws4.Cells.Item($r, $ws4PrevCol).NumberFormat
}

# ---------------------------------------------------------------------
# Sheet 5: "Lives Lost by Ward" — new column Z (26), header row = 2.
# Previous column is Y (25).
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Lives Lost by Ward")
$ws5NewCol = 26
$ws5PrevCol = 25

$sheet5Data = @{
  2 = 43964
  3 = 358
  4 = 33
  5 = 18
  6 = 23
  7 = 48
  8 = 56
  9 = 38
  10 = 47
  11 = 80
  12 = 15
  13 = 0
}

foreach ($r in $sheet5Data.Keys) {
  $cell = $ws5.Cells.Item($r, $ws5NewCol)
  $cell.Value = $sheet5Data[$r]
  $cell.NumberFormat = $ws5.Cells.Item($r, $ws5PrevCol).NumberFormat
}
